$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.897.63"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.632.25"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'211.38"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'23.42"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").Value = "'0.0883"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "1.864.37"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "1.639.42"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "27.908.43"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "'229.02"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("D20").Value = "0.0₃0719"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").Value = "'10.05"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").Value = "'154.96"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'15.52"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "'3.11"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("D34").Value = "1.392.75"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("E36").Value = "  +10.09%  "
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("D43").Value = "'65.71"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").Value = "1.774.51"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("E47").Value = "  -3.13%  "
$ws.Range("D48").Value = "'88.73"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("E51").Value = "  +0.87%  "
